$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'305.87"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'1.25%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'36.17"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'-1.16%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'5.068"
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').Value = "'0.07974"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'2.86%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'2.170"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'3.18%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'8.030"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'1.55%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D8').Value = "'4.163"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'3.07%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D9').Value = "'0.9301"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'1.01%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').Value = "'0.09851"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'0.91%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D11').Value = "'0.1868"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'0.26%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D12').Value = "'0.09073"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'5.80%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = "'0.03620"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'3.09%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = "'0.09922"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'-0.24%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D15').Value = "'0.001434"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'-2.32%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').Value = "'0.005664"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'0.34%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').Value = "'3.455"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'-0.17%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'8.86%"
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'0.3371"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'-1.08%"
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.1358"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'1.18%"
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'5.088"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'6.33%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.2194"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'-0.18%"
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'0.04580"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'-0.10%"
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'0.001242"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'1.15%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.004754"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'-6.48%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'-6.68%"
$ws.Range('E26').Style = 'Normal'
$ws.Range('D39').Value = "'0.01946"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'10.06%"
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.04907"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'5.47%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'0.007817"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'4.52%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.1394"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'0.11%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'0.007820"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'1.52%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.002102"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'-5.92%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.01144"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'9.77%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.00006217"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'0.83%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.00000000752"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'0.50%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'52.15"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'36.56%"
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.001804"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'-9.55%"
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.00002105"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'0.50%"
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.0002004"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'0.50%"
$ws.Range('E51').Style = 'Normal'
